# Updated symbol list on Sun Dec 18 21:38:10 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $value)
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Simple price-only updates (column D) ---
Set-TextCell "D2" "248.47"
Set-TextCell "D3" "21.74"
Set-TextCell "D4" "5.578"
Set-TextCell "D5" "0.05669"
Set-TextCell "D6" "6.433"
Set-TextCell "D7" "0.8010"
Set-TextCell "D9" "0.1435"
Set-TextCell "D10" "0.07268"
Set-TextCell "D11" "0.03129"
Set-TextCell "D12" "0.02923"
Set-TextCell "D13" "0.09277"
Set-TextCell "D14" "0.001676"
Set-TextCell "D15" "3.210"
Set-TextCell "D16" "0.04750"
Set-TextCell "D17" "0.0005828"
Set-TextCell "D18" "0.006462"
Set-TextCell "D19" "0.005058"
Set-TextCell "D20" "0.001056"
Set-TextCell "D21" "0.0001502"
Set-TextCell "D22" "3.977"
Set-TextCell "D23" "3.380"
Set-TextCell "D24" "2.114"
Set-TextCell "D27" "0.0003205"
Set-TextCell "D40" "0.04127"
Set-TextCell "D41" "0.006929"

# --- Row 42/43 swap: CEJI and BKEXToken change places in ranking ---
Set-TextCell "B42" "BKEXToken"
Set-TextCell "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D42" "0.1045"
Set-TextCell "E42" "41BKEXTokenBKK"

Set-TextCell "B43" "CEJI"
Set-TextCell "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell "D43" "0.002974"
Set-TextCell "E43" "42CEJICEJI"

# --- Remaining price-only updates ---
Set-TextCell "D44" "0.009394"
Set-TextCell "D45" "0.00005657"

# Row 47: price update and the "Bestin24h" marker moves here
Set-TextCell "D47" "0.7863"
Set-TextCell "E47" "46CoinbaseStockTokenCOINBestin24h"

Set-TextCell "D48" "0.01670"
Set-TextCell "D49" "0.00002103"
